# Apply cell-level updates to match the refreshed cryptos data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.790.14"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.634.99"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.03%  "
$r = $ws.Range("D5")
$r.Value = "'215.34"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -0.05%  "
$r = $ws.Range("D8")
$r.Value = "'0.257"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$r = $ws.Range("D9")
$r.Value = "'0.0642"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$r = $ws.Range("D10")
$r.Value = "'19.84"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +1.55%  "
$r = $ws.Range("D11")
$r.Value = "'0.0786"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.640.37"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$r = $ws.Range("D13")
$r.Value = "'4.25"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "1.860.52"
$ws.Range("E14").Value = "  -0.13%  "
$r = $ws.Range("D15")
$r.Value = "'0.555"
$r.Style = "Normal"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "0.0₃0775"
$ws.Range("E16").Value = "  +2.15%  "
$r = $ws.Range("D17")
$r.Value = "'63.09"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "25.803.49"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  +2.75%  "
$r = $ws.Range("D21")
$r.Value = "'194.57"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -1.60%  "
$r = $ws.Range("D26")
$r.Value = "'139.70"
$r.Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -5.05%  "
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("E30").Value = "  -0.07%  "
$r = $ws.Range("D31")
$r.Value = "'0.0494"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("E35").Value = "  +0.71%  "
$r = $ws.Range("D36")
$r.Value = "'0.897"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("E37").Value = "  -0.33%  "
$r = $ws.Range("D38")
$r.Value = "'0.551"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").Value = "1.110.89"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("E42").Value = "  +1.08%  "
$r = $ws.Range("D43")
$r.Value = "'99.26"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  -1.93%  "
$r = $ws.Range("D46")
$r.Value = "'55.56"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "
$r = $ws.Range("D47")
$r.Value = "'2.52"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +13.73%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$r = $ws.Range("D48")
$r.Value = "'7.72"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$r = $ws.Range("D49")
$r.Value = "'0.418"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -5.13%  "
$r = $ws.Range("D50")
$r.Value = "'0.0503"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -0.19%  "
